$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells: force text storage so numeric-looking strings
# (e.g. "592.72", "14.38") are not auto-converted to Excel numbers,
# matching the source workbook where these are inline/shared strings.
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '64.320.70'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.151.85'
$ws.Range('D3').Style = "Normal"
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '592.72'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.99'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.35%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '3.140.04'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.79%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.530'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  +0.46%  '
$ws.Range('E11').Value = '  +3.62%  '
$ws.Range('E12').Value = '  -1.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000247'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.58%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '37.46'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.672.90'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.95%  '
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '64.083.24'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.35%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.144.32'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '470.01'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.38'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.37%  '
$ws.Range('E22').Value = '  +0.60%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.64'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.30%  '
$ws.Range('E24').Value = '  +9.20%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '13.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.09%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '81.48'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('E27').Value = '  +0.02%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.96'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +11.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.46'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +8.28%  '
$ws.Range('E30').Value = '  +1.05%  '
$ws.Range('E31').Value = '  +0.74%  '
$ws.Range('E32').Value = '  +0.16%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.69'
$ws.Range('D33').Style = "Normal"
$ws.Range('E34').Value = '  +0.94%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0₃0850'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -1.93%  '
$ws.Range('E36').Value = '  +1.64%  '
$ws.Range('E37').Value = '  +2.04%  '
$ws.Range('E38').Value = '  -2.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.21'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '51.42'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.86%  '
$ws.Range('B41').Value = 'Cosmos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.30'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +6.96%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '456.45'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.293'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +6.33%  '
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.928.55'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.67%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.59'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +13.59%  '
$ws.Range('E47').Value = '  -2.49%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '134.23'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +8.79%  '
$ws.Range('E50').Value = '  +3.17%  '
$ws.Range('E51').Value = '  -0.11%  '
